$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phase 1: introduce new shared strings in the exact order used by the source workbook
$intro = @(
    ,@(52, 4, "I appreciate hugs or cuddles more than long and sweet messages")
    ,@(52, 3, "LoveTest")
    ,@(52, 6, "Physical Touch")
    ,@(53, 6, "Words of Affirmation")
    ,@(54, 4, "I appreciate holding hands than being complimented by the person I like")
    ,@(56, 4, "I would rather stay at home and cuddle with my partner than sending motivational messages to my partner")
    ,@(58, 4, "I like putting my arms around my partner more than encouraging her verbally")
    ,@(60, 4, "I would rather go on a hike with my partner than helping him/her with his/her project")
    ,@(62, 4, "I prefer having heart-to-heart talks with my partner to helping him/her for his/her exams")
    ,@(64, 4, "I prefer being alone with my partner to helping him/her do the dishes")
    ,@(66, 4, "I appreciate having dinner alone with my partner more than fetching him/her from work/school")
    ,@(68, 4, "I would rather hug my partner than spend alone time with my partner")
    ,@(70, 4, "I would rather encourage my partner when my partner is stressed than help her out")
    ,@(60, 6, "Quality Time")
    ,@(61, 6, "Act of Service")
    ,@(72, 3, "Job")
    ,@(72, 4, "I prefer working alone than in groups.")
    ,@(74, 4, "I enjoy meeting clients.")
    ,@(76, 4, "I prefer my supervisors to be honest regarding my performances.")
    ,@(78, 4, "I prefer to work in a traditional company than a start up.")
    ,@(80, 4, "I enjoy attaining physical meetings and conference rather than virtual ones.")
    ,@(82, 4, "I prefer my colleagues to not interfere with my personal life.")
    ,@(84, 4, "I prefer companies that allows me to share personal opinions and gives room for progression.")
    ,@(86, 4, "I prefer working in a community rather than a hierarchical environment.")
    ,@(88, 4, "I think everyone is equal without having my boss to tell me so.")
    ,@(90, 4, "I do not like my boss to point out my flaws infront of everyone.")
    ,@(72, 6, "Outgoing")
    ,@(73, 6, "Asocial")
    ,@(76, 6, "Idealistic")
    ,@(77, 6, "Realistic")
)
foreach ($i in $intro) {
    $ws.Cells.Item($i[0], $i[1]).Value = $i[2]
}

# Phase 2: write full row data (values already introduced will reuse shared strings)
$rows = @(
    ,@(52, 51, 1, "LoveTest", "I appreciate hugs or cuddles more than long and sweet messages", 1, "Physical Touch")
    ,@(53, 52, 1, "LoveTest", "I appreciate hugs or cuddles more than long and sweet messages", 0, "Words of Affirmation")
    ,@(54, 53, 2, "LoveTest", "I appreciate holding hands than being complimented by the person I like", 1, "Physical Touch")
    ,@(55, 54, 2, "LoveTest", "I appreciate holding hands than being complimented by the person I like", 0, "Words of Affirmation")
    ,@(56, 55, 3, "LoveTest", "I would rather stay at home and cuddle with my partner than sending motivational messages to my partner", 1, "Physical Touch")
    ,@(57, 56, 3, "LoveTest", "I would rather stay at home and cuddle with my partner than sending motivational messages to my partner", 0, "Words of Affirmation")
    ,@(58, 57, 4, "LoveTest", "I like putting my arms around my partner more than encouraging her verbally", 1, "Physical Touch")
    ,@(59, 58, 4, "LoveTest", "I like putting my arms around my partner more than encouraging her verbally", 0, "Words of Affirmation")
    ,@(60, 59, 5, "LoveTest", "I would rather go on a hike with my partner than helping him/her with his/her project", 1, "Quality Time")
    ,@(61, 60, 5, "LoveTest", "I would rather go on a hike with my partner than helping him/her with his/her project", 0, "Act of Service")
    ,@(62, 61, 6, "LoveTest", "I prefer having heart-to-heart talks with my partner to helping him/her for his/her exams", 1, "Quality Time")
    ,@(63, 62, 6, "LoveTest", "I prefer having heart-to-heart talks with my partner to helping him/her for his/her exams", 0, "Act of Service")
    ,@(64, 63, 7, "LoveTest", "I prefer being alone with my partner to helping him/her do the dishes", 1, "Quality Time")
    ,@(65, 64, 7, "LoveTest", "I prefer being alone with my partner to helping him/her do the dishes", 0, "Act of Service")
    ,@(66, 65, 8, "LoveTest", "I appreciate having dinner alone with my partner more than fetching him/her from work/school", 1, "Quality Time")
    ,@(67, 66, 8, "LoveTest", "I appreciate having dinner alone with my partner more than fetching him/her from work/school", 0, "Act of Service")
    ,@(68, 67, 9, "LoveTest", "I would rather hug my partner than spend alone time with my partner", 1, "Physical Touch")
    ,@(69, 68, 9, "LoveTest", "I would rather hug my partner than spend alone time with my partner", 0, "Quality Time")
    ,@(70, 69, 10, "LoveTest", "I would rather encourage my partner when my partner is stressed than help her out", 1, "Words of Affirmation")
    ,@(71, 70, 10, "LoveTest", "I would rather encourage my partner when my partner is stressed than help her out", 0, "Act of Service")
    ,@(72, 71, 1, "Job", "I prefer working alone than in groups.", 1, "Outgoing")
    ,@(73, 72, 1, "Job", "I prefer working alone than in groups.", 0, "Asocial")
    ,@(74, 73, 2, "Job", "I enjoy meeting clients.", 1, "Outgoing")
    ,@(75, 74, 2, "Job", "I enjoy meeting clients.", 0, "Asocial")
    ,@(76, 75, 3, "Job", "I prefer my supervisors to be honest regarding my performances.", 1, "Idealistic")
    ,@(77, 76, 3, "Job", "I prefer my supervisors to be honest regarding my performances.", 0, "Realistic")
    ,@(78, 77, 4, "Job", "I prefer to work in a traditional company than a start up.", 1, "Idealistic")
    ,@(79, 78, 4, "Job", "I prefer to work in a traditional company than a start up.", 0, "Realistic")
    ,@(80, 79, 5, "Job", "I enjoy attaining physical meetings and conference rather than virtual ones.", 1, "Outgoing")
    ,@(81, 80, 5, "Job", "I enjoy attaining physical meetings and conference rather than virtual ones.", 0, "Asocial")
    ,@(82, 81, 6, "Job", "I prefer my colleagues to not interfere with my personal life.", 1, "Outgoing")
    ,@(83, 82, 6, "Job", "I prefer my colleagues to not interfere with my personal life.", 0, "Asocial")
    ,@(84, 83, 7, "Job", "I prefer companies that allows me to share personal opinions and gives room for progression.", 1, "Idealistic")
    ,@(85, 84, 7, "Job", "I prefer companies that allows me to share personal opinions and gives room for progression.", 0, "Realistic")
    ,@(86, 85, 8, "Job", "I prefer working in a community rather than a hierarchical environment.", 1, "Idealistic")
    ,@(87, 86, 8, "Job", "I prefer working in a community rather than a hierarchical environment.", 0, "Realistic")
    ,@(88, 87, 9, "Job", "I think everyone is equal without having my boss to tell me so.", 1, "Idealistic")
    ,@(89, 88, 9, "Job", "I think everyone is equal without having my boss to tell me so.", 0, "Outgoing")
    ,@(90, 89, 10, "Job", "I do not like my boss to point out my flaws infront of everyone.", 1, "Asocial")
    ,@(91, 90, 10, "Job", "I do not like my boss to point out my flaws infront of everyone.", 0, "Realistic")
)
foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Range($ws.Cells.Item($rowNum,3), $ws.Cells.Item($rowNum,3)).Font.Name = "Arial"
    $ws.Range($ws.Cells.Item($rowNum,3), $ws.Cells.Item($rowNum,3)).Font.Size = 10
    $ws.Range($ws.Cells.Item($rowNum,4), $ws.Cells.Item($rowNum,4)).Font.Name = "Arial"
    $ws.Range($ws.Cells.Item($rowNum,4), $ws.Cells.Item($rowNum,4)).Font.Size = 10
    $ws.Range($ws.Cells.Item($rowNum,6), $ws.Cells.Item($rowNum,6)).Font.Name = "Arial"
    $ws.Range($ws.Cells.Item($rowNum,6), $ws.Cells.Item($rowNum,6)).Font.Size = 10
}

$ws.Range("G90").Select()
Write-Host "Done"